# Update the "Model holdings provided as of" date in the confidential disclosure text (cell A59)
# and refresh the Weight (D) and Percent Change (E) columns for all holdings (rows 2-56).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The worksheet is protected; unprotect temporarily so the cell values can be updated.
$ws.Unprotect()

# --- Update confidential disclosure date string (A59) ---
$oldText = $ws.Range("A59").Value2
$newText = $oldText -replace "2021-06-14", "2021-07-07"
$ws.Range("A59").Value = $newText

# --- Update Weight (column D) and Percent Change (column E) for each holding row ---
    $ws.Cells.Item(2, 4).Value = 0.01426826839319593
    $ws.Cells.Item(2, 5).Value = 0
    $ws.Cells.Item(3, 4).Value = 0.05591746562479984
    $ws.Cells.Item(3, 5).Value = 0
    $ws.Cells.Item(4, 4).Value = 0.01573162670698091
    $ws.Cells.Item(4, 5).Value = 0
    $ws.Cells.Item(5, 4).Value = 0.008876841001099971
    $ws.Cells.Item(5, 5).Value = 0
    $ws.Cells.Item(6, 4).Value = 0.01542682758047619
    $ws.Cells.Item(6, 5).Value = 0
    $ws.Cells.Item(7, 4).Value = 0.01944489424475253
    $ws.Cells.Item(7, 5).Value = 0
    $ws.Cells.Item(8, 4).Value = 0.004526592577672236
    $ws.Cells.Item(8, 5).Value = 0
    $ws.Cells.Item(9, 4).Value = 0.006721888662102477
    $ws.Cells.Item(9, 5).Value = 0
    $ws.Cells.Item(10, 4).Value = 0.01364093227926092
    $ws.Cells.Item(10, 5).Value = 0
    $ws.Cells.Item(11, 4).Value = 0.008192883687881135
    $ws.Cells.Item(11, 5).Value = 0
    $ws.Cells.Item(12, 4).Value = 0.01362559800497879
    $ws.Cells.Item(12, 5).Value = 0
    $ws.Cells.Item(13, 4).Value = 0.002664330156520169
    $ws.Cells.Item(13, 5).Value = 0
    $ws.Cells.Item(14, 4).Value = 0.005599139873294584
    $ws.Cells.Item(14, 5).Value = 0
    $ws.Cells.Item(15, 4).Value = 0.01353581339504512
    $ws.Cells.Item(15, 5).Value = 0
    $ws.Cells.Item(16, 4).Value = 0.009338846864143921
    $ws.Cells.Item(16, 5).Value = 0
    $ws.Cells.Item(17, 4).Value = 0.0222976777641765
    $ws.Cells.Item(17, 5).Value = 0
    $ws.Cells.Item(18, 4).Value = 0.008516911507533299
    $ws.Cells.Item(18, 5).Value = 0
    $ws.Cells.Item(19, 4).Value = 0.01635363842012352
    $ws.Cells.Item(19, 5).Value = 0
    $ws.Cells.Item(20, 4).Value = 0.01432537639483791
    $ws.Cells.Item(20, 5).Value = 0
    $ws.Cells.Item(21, 4).Value = 0.006533739550672528
    $ws.Cells.Item(21, 5).Value = 0
    $ws.Cells.Item(22, 4).Value = 0.01383151540248168
    $ws.Cells.Item(22, 5).Value = 0
    $ws.Cells.Item(23, 4).Value = 0.01829351539225518
    $ws.Cells.Item(23, 5).Value = 0
    $ws.Cells.Item(24, 4).Value = 0.008831233705169587
    $ws.Cells.Item(24, 5).Value = 0
    $ws.Cells.Item(25, 4).Value = 0.02079996945899365
    $ws.Cells.Item(25, 5).Value = 0
    $ws.Cells.Item(26, 4).Value = 0.01349991172113061
    $ws.Cells.Item(26, 5).Value = 0
    $ws.Cells.Item(27, 4).Value = 0.02081353907472744
    $ws.Cells.Item(27, 5).Value = 0
    $ws.Cells.Item(28, 4).Value = 0.06049371204300471
    $ws.Cells.Item(28, 5).Value = 0
    $ws.Cells.Item(29, 4).Value = 0.01856974530536523
    $ws.Cells.Item(29, 5).Value = 0
    $ws.Cells.Item(30, 4).Value = 0.03160807711551528
    $ws.Cells.Item(30, 5).Value = 0
    $ws.Cells.Item(31, 4).Value = 0.01589418826945993
    $ws.Cells.Item(31, 5).Value = 0
    $ws.Cells.Item(32, 4).Value = 0.01338018876616993
    $ws.Cells.Item(32, 5).Value = 0
    $ws.Cells.Item(33, 4).Value = 0.01776147083993055
    $ws.Cells.Item(33, 5).Value = 0
    $ws.Cells.Item(34, 4).Value = 0.04607219218242962
    $ws.Cells.Item(34, 5).Value = 0
    $ws.Cells.Item(35, 4).Value = 0.01028491682208608
    $ws.Cells.Item(35, 5).Value = 0
    $ws.Cells.Item(36, 4).Value = 0.009859086459283106
    $ws.Cells.Item(36, 5).Value = 0
    $ws.Cells.Item(37, 4).Value = 0.009866875297013711
    $ws.Cells.Item(37, 5).Value = 0
    $ws.Cells.Item(38, 4).Value = 0.006665845540618977
    $ws.Cells.Item(38, 5).Value = 0
    $ws.Cells.Item(39, 4).Value = 0.01133796159823453
    $ws.Cells.Item(39, 5).Value = 0
    $ws.Cells.Item(40, 4).Value = 0.01708639766974819
    $ws.Cells.Item(40, 5).Value = 0
    $ws.Cells.Item(41, 4).Value = 0.01730764934153321
    $ws.Cells.Item(41, 5).Value = 0
    $ws.Cells.Item(42, 4).Value = 0.03647366668535319
    $ws.Cells.Item(42, 5).Value = 0
    $ws.Cells.Item(43, 4).Value = 0.01118900798717501
    $ws.Cells.Item(43, 5).Value = 0
    $ws.Cells.Item(44, 4).Value = 0.02259858247181593
    $ws.Cells.Item(44, 5).Value = 0
    $ws.Cells.Item(45, 4).Value = 0.01295965237901198
    $ws.Cells.Item(45, 5).Value = 0
    $ws.Cells.Item(46, 4).Value = 0.00806719740403296
    $ws.Cells.Item(46, 5).Value = 0
    $ws.Cells.Item(47, 4).Value = 0.01169347944543035
    $ws.Cells.Item(47, 5).Value = 0
    $ws.Cells.Item(48, 4).Value = 0.009288462820074063
    $ws.Cells.Item(48, 5).Value = 0
    $ws.Cells.Item(49, 4).Value = 0.01629595234068122
    $ws.Cells.Item(49, 5).Value = 0
    $ws.Cells.Item(50, 4).Value = 0.008108575604476803
    $ws.Cells.Item(50, 5).Value = 0
    $ws.Cells.Item(51, 4).Value = 0.0105094544097887
    $ws.Cells.Item(51, 5).Value = 0
    $ws.Cells.Item(52, 4).Value = 0.00834516155054396
    $ws.Cells.Item(52, 5).Value = 0
    $ws.Cells.Item(53, 4).Value = 0.009660379821710498
    $ws.Cells.Item(53, 5).Value = 0
    $ws.Cells.Item(54, 4).Value = 0.1320557276029691
    $ws.Cells.Item(54, 5).Value = 0
    $ws.Cells.Item(55, 4).Value = 0.04495741478223664
    $ws.Cells.Item(55, 5).Value = 0
    $ws.Cells.Item(56, 4).Value = 1
    $ws.Cells.Item(56, 5).Value = 0


# Re-apply sheet protection to restore the worksheet to its protected state.
$ws.Protect()
